$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three lines that used to live one-per-row in A6:A8 are regrouped into
# two rows (6 and 7), using column B for the second half of each line.
# Row 6: A6 keeps "EDM is a genre..."; B6 picks up what used to be A7's text
# (and takes on A6's formatting, since it's now a "title" style cell).
# Row 7: A7 picks up what used to be A8's text; B7 gets a small blank pad.
# Row 8's A cell is then emptied out completely (row stays, cell goes).
$a7Text = $ws.Range("A7").Value()
$a8Text = $ws.Range("A8").Value()

$ws.Range("A6").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = $a7Text

$ws.Range("A7").Value = $a8Text
$ws.Range("B7").Value = "   "
$ws.Range("A8").Clear()

# Rows 6 and 7 now render at the shorter auto height (matches rows 9-11).
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8

# Column widths nudged slightly (A a bit wider, B a touch narrower).
$ws.Columns.Item(1).ColumnWidth = 77.43707482993197
$ws.Columns.Item(2).ColumnWidth = 13.069727891156466

# Selection moved to B7.
$ws.Range("B7").Select()
